$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.51%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.38%"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.15%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07837"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.11%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.844"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-15.28%"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.076"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.85%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.804"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.65%"

$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.840"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "7.93%"

$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9249"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.46%"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1072"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "8.46%"

$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1851"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.78%"

$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09421"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.78%"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03577"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.41%"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09921"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.04%"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001415"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.63%"

$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005741"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.35%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.459"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.18%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.74%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1295"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.56%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.111"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.87%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2200"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.45%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04567"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.49%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001228"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.07%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004655"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.00%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001255"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.56%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004465"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.94%"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.03%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04705"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.89%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007605"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009994"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "27.93%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1334"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.37%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002122"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.77%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01126"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.53%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006186"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.05%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.22%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "66.41"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "27.70%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001305"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-27.62%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002108"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.22%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.22%"
